$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 335 (existing rows 335-338 shift down to 337-340)
$ws.Rows.Item(335).Insert()
$ws.Rows.Item(335).Insert()

# New row 335 data
$ws.Cells.Item(335, 1).Value = 6
$ws.Cells.Item(335, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(335, 3).Value = "Metropolitana"
$ws.Cells.Item(335, 4).Value = 44448
$ws.Cells.Item(335, 5).Value = 13
$ws.Cells.Item(335, 6).Value = 100112012
$ws.Cells.Item(335, 7).Value = "Espinaca"
$ws.Cells.Item(335, 8).Value = "Sin especificar"
$ws.Cells.Item(335, 9).Value = "Primera"
$ws.Cells.Item(335, 10).Value = 440
$ws.Cells.Item(335, 11).Value = 4000
$ws.Cells.Item(335, 12).Value = 4500
$ws.Cells.Item(335, 13).Value = 4205
$ws.Cells.Item(335, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(335, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(335, 16).Value = 420
$ws.Cells.Item(335, 17).Value = 10
$ws.Cells.Item(335, 18).Value = "Hortaliza"

# New row 336 data
$ws.Cells.Item(336, 1).Value = 6
$ws.Cells.Item(336, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(336, 3).Value = "Metropolitana"
$ws.Cells.Item(336, 4).Value = 44448
$ws.Cells.Item(336, 5).Value = 13
$ws.Cells.Item(336, 6).Value = 100112012
$ws.Cells.Item(336, 7).Value = "Espinaca"
$ws.Cells.Item(336, 8).Value = "Sin especificar"
$ws.Cells.Item(336, 9).Value = "Primera"
$ws.Cells.Item(336, 10).Value = 510
$ws.Cells.Item(336, 11).Value = 4000
$ws.Cells.Item(336, 12).Value = 4500
$ws.Cells.Item(336, 13).Value = 4216
$ws.Cells.Item(336, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(336, 15).Value = "Región Metropolitana"
$ws.Cells.Item(336, 16).Value = 422
$ws.Cells.Item(336, 17).Value = 10
$ws.Cells.Item(336, 18).Value = "Hortaliza"

Write-Output "done"
